$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# G3 and G5 hold the "Param Value" for the two "Global / DoSleep / millis"
# rows. Change the sleep duration from 2000 to 1000, keeping the cell a
# text value (as it was originally) rather than letting it be reinterpreted
# as a number.
$ws.Range("G3").Value = "'1000"
$ws.Range("G3").Style = "Normal"

$ws.Range("G5").Value = "'1000"
$ws.Range("G5").Style = "Normal"
